# Auto-generated Excel COM-interop script
# Fixes formatting introduced by a scraping bug that applied a blanket
# str.Replace(".", "").Replace(",", ".") over shared strings:
#   - Column H ("Importe") values move from Spanish-style text
#     ("1.234,56") to plain decimal text ("1234.56"), remaining TEXT.
#   - A few "Razon social" (column E) strings had commas mistakenly
#     turned into periods by the same bug; restore that mapping here
#     (comma -> period, and any other "." in the string removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("Importe"): keep values as TEXT while rewriting their format ---
# Mark the whole range as Text first so Excel does not silently coerce the
# numeric-looking strings ("1500.00") into real numbers (which would drop
# the trailing ".00" and change the stored cell type).
$importeRange = $ws.Range("H2:H269")
$importeRange.NumberFormat = "@"

$importeChanges = @(
    ,@('H2', '1500.00')
    ,@('H3', '3005.00')
    ,@('H4', '3010.00')
    ,@('H5', '1200.00')
    ,@('H6', '850000.00')
    ,@('H7', '486578.83')
    ,@('H8', '9438.00')
    ,@('H9', '19600.00')
    ,@('H10', '1300.00')
    ,@('H11', '465298.50')
    ,@('H12', '194123.49')
    ,@('H13', '18636.80')
    ,@('H14', '123941.00')
    ,@('H15', '130.00')
    ,@('H16', '3500.00')
    ,@('H17', '45386.86')
    ,@('H18', '489.00')
    ,@('H19', '38542.73')
    ,@('H20', '20237.67')
    ,@('H21', '7900.00')
    ,@('H22', '960.00')
    ,@('H23', '4950.00')
    ,@('H24', '9600.00')
    ,@('H25', '1364.42')
    ,@('H26', '116.00')
    ,@('H27', '52.00')
    ,@('H28', '7045.00')
    ,@('H29', '3605.55')
    ,@('H30', '19.61')
    ,@('H31', '365.00')
    ,@('H32', '3130.00')
    ,@('H33', '527600.00')
    ,@('H34', '19085.00')
    ,@('H35', '80.00')
    ,@('H36', '1573.00')
    ,@('H37', '8887.78')
    ,@('H38', '980.00')
    ,@('H39', '3449.40')
    ,@('H40', '35305.70')
    ,@('H41', '2831.09')
    ,@('H42', '257.00')
    ,@('H43', '3800.00')
    ,@('H44', '883.80')
    ,@('H45', '6300.00')
    ,@('H46', '9606.00')
    ,@('H47', '6603.30')
    ,@('H48', '150.00')
    ,@('H49', '40.00')
    ,@('H50', '4640.90')
    ,@('H51', '10609.33')
    ,@('H52', '2089.00')
    ,@('H53', '160.00')
    ,@('H54', '2115.00')
    ,@('H55', '1220.00')
    ,@('H56', '11854.02')
    ,@('H57', '800.00')
    ,@('H58', '6294.99')
    ,@('H59', '32445.00')
    ,@('H60', '19750.00')
    ,@('H61', '6800.00')
    ,@('H62', '1000.00')
    ,@('H63', '114.59')
    ,@('H64', '696.00')
    ,@('H65', '2189.39')
    ,@('H66', '391.00')
    ,@('H67', '1200.00')
    ,@('H68', '4000.00')
    ,@('H69', '4101.00')
    ,@('H70', '108000.00')
    ,@('H71', '690.00')
    ,@('H72', '889.80')
    ,@('H73', '1053.60')
    ,@('H74', '1800.00')
    ,@('H75', '39796.00')
    ,@('H76', '5800.00')
    ,@('H77', '32380.00')
    ,@('H78', '48800.00')
    ,@('H79', '2812.10')
    ,@('H80', '13770.00')
    ,@('H81', '20625.00')
    ,@('H82', '21943.00')
    ,@('H83', '3990.00')
    ,@('H84', '81407.00')
    ,@('H85', '61257.79')
    ,@('H86', '572.03')
    ,@('H87', '45900.00')
    ,@('H88', '28314.00')
    ,@('H89', '2670.00')
    ,@('H90', '21.17')
    ,@('H91', '2939.00')
    ,@('H92', '600235.87')
    ,@('H93', '3189.89')
    ,@('H94', '247.76')
    ,@('H95', '18.12')
    ,@('H96', '23953.63')
    ,@('H97', '660.00')
    ,@('H98', '900.00')
    ,@('H99', '8196.00')
    ,@('H100', '292.74')
    ,@('H101', '14835.92')
    ,@('H102', '55520.80')
    ,@('H103', '2700.00')
    ,@('H104', '246.00')
    ,@('H105', '800.00')
    ,@('H106', '3400.00')
    ,@('H107', '200.00')
    ,@('H108', '346.00')
    ,@('H109', '113151.40')
    ,@('H110', '6978.00')
    ,@('H111', '6880.00')
    ,@('H112', '49202.35')
    ,@('H113', '10061.00')
    ,@('H114', '10144.80')
    ,@('H115', '10739.81')
    ,@('H116', '1035.20')
    ,@('H117', '11246.22')
    ,@('H118', '2750.00')
    ,@('H119', '354.00')
    ,@('H120', '7753.80')
    ,@('H121', '6390.00')
    ,@('H122', '4686.00')
    ,@('H123', '9300.00')
    ,@('H124', '63000.00')
    ,@('H125', '641152.00')
    ,@('H126', '84700.00')
    ,@('H127', '11609.00')
    ,@('H128', '4000.00')
    ,@('H129', '16371.00')
    ,@('H130', '11902.00')
    ,@('H131', '5700.00')
    ,@('H132', '15372.00')
    ,@('H133', '66400.00')
    ,@('H134', '3500.00')
    ,@('H135', '8000.00')
    ,@('H136', '2190.00')
    ,@('H137', '2500.00')
    ,@('H138', '36500.00')
    ,@('H139', '8000.00')
    ,@('H140', '3250.00')
    ,@('H141', '261940.00')
    ,@('H142', '16626.05')
    ,@('H143', '309.00')
    ,@('H144', '11735.55')
    ,@('H145', '6748.88')
    ,@('H146', '3360.00')
    ,@('H147', '700.00')
    ,@('H148', '2020.00')
    ,@('H149', '40000.00')
    ,@('H150', '10000.00')
    ,@('H151', '6000.00')
    ,@('H152', '14000.00')
    ,@('H153', '6000.00')
    ,@('H154', '25116.58')
    ,@('H155', '4500.00')
    ,@('H156', '3000.00')
    ,@('H157', '13260.00')
    ,@('H158', '8000.00')
    ,@('H159', '3500.00')
    ,@('H160', '3000.00')
    ,@('H161', '4000.00')
    ,@('H162', '4500.00')
    ,@('H163', '3500.00')
    ,@('H164', '3500.00')
    ,@('H165', '4500.00')
    ,@('H166', '5000.00')
    ,@('H167', '10500.00')
    ,@('H168', '4000.00')
    ,@('H169', '6000.00')
    ,@('H170', '4500.00')
    ,@('H171', '2000.00')
    ,@('H172', '15000.00')
    ,@('H173', '4000.00')
    ,@('H174', '4000.00')
    ,@('H175', '2500.00')
    ,@('H176', '4000.00')
    ,@('H177', '2359.50')
    ,@('H178', '4750.00')
    ,@('H179', '1080.00')
    ,@('H180', '144.42')
    ,@('H181', '9515.00')
    ,@('H182', '1460.00')
    ,@('H183', '4997.60')
    ,@('H184', '500.00')
    ,@('H185', '19830.00')
    ,@('H186', '22.34')
    ,@('H187', '11700.00')
    ,@('H188', '1078.00')
    ,@('H189', '5155.00')
    ,@('H190', '6433.85')
    ,@('H191', '65.32')
    ,@('H192', '26897.31')
    ,@('H193', '6782.14')
    ,@('H194', '3690.50')
    ,@('H195', '8240.00')
    ,@('H196', '5200.00')
    ,@('H197', '13350.00')
    ,@('H198', '1215.89')
    ,@('H199', '4106.53')
    ,@('H200', '78.00')
    ,@('H201', '3551.20')
    ,@('H202', '48249.25')
    ,@('H203', '630.00')
    ,@('H204', '50207.10')
    ,@('H205', '2135.98')
    ,@('H206', '36630.00')
    ,@('H207', '65894.00')
    ,@('H208', '240.00')
    ,@('H209', '8844.00')
    ,@('H210', '62594.13')
    ,@('H211', '8500.00')
    ,@('H212', '30000.00')
    ,@('H213', '30000.00')
    ,@('H214', '30000.00')
    ,@('H215', '30000.00')
    ,@('H216', '30000.00')
    ,@('H217', '30000.00')
    ,@('H218', '60000.00')
    ,@('H219', '60000.00')
    ,@('H220', '60000.00')
    ,@('H221', '30000.00')
    ,@('H222', '44400.00')
    ,@('H223', '3019313.66')
    ,@('H224', '7200.00')
    ,@('H225', '18400.00')
    ,@('H226', '139500.00')
    ,@('H227', '4800.00')
    ,@('H228', '135500.00')
    ,@('H229', '146000.00')
    ,@('H230', '135500.00')
    ,@('H231', '141100.00')
    ,@('H232', '135500.00')
    ,@('H233', '135500.00')
    ,@('H234', '239000.00')
    ,@('H235', '239000.00')
    ,@('H236', '347000.00')
    ,@('H237', '135500.00')
    ,@('H238', '135500.00')
    ,@('H239', '135500.00')
    ,@('H240', '135500.00')
    ,@('H241', '135500.00')
    ,@('H242', '239000.00')
    ,@('H243', '342500.00')
    ,@('H244', '239000.00')
    ,@('H245', '135500.00')
    ,@('H246', '244000.00')
    ,@('H247', '135500.00')
    ,@('H248', '135500.00')
    ,@('H249', '140250.00')
    ,@('H250', '135500.00')
    ,@('H251', '34623.36')
    ,@('H252', '25000.00')
    ,@('H253', '69450.00')
    ,@('H254', '16000.00')
    ,@('H255', '20000.00')
    ,@('H256', '83734.00')
    ,@('H257', '17000.00')
    ,@('H258', '1500.00')
    ,@('H259', '30000.00')
    ,@('H260', '1184830.00')
    ,@('H261', '2750.00')
    ,@('H262', '7405.00')
    ,@('H263', '4500.00')
    ,@('H264', '25878.00')
    ,@('H265', '25000.00')
    ,@('H266', '18000.00')
    ,@('H267', '3570.00')
    ,@('H268', '47845.00')
    ,@('H269', '23500.00')
)

foreach ($pair in $importeChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Restore the default (unstyled) look - only the underlying value/type matters,
# the visible number format was never meant to change.
$importeRange.Style = "Normal"

# --- Column E ("Razon social"): comma -> period corrections ---
$razonSocialChanges = @(
    ,@('E31', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO')
    ,@('E83', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO')
    ,@('E262', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO')
    ,@('E81', 'FERNANDEZ. MARIO HUGO')
    ,@('E173', 'RICCOTTI. MARIANA EDITH')
    ,@('E181', 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN')
    ,@('E195', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH')
)

foreach ($pair in $razonSocialChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

